$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Reln"
$ws.Range("C2").Value = "Lrp8"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.03565433333333334
$ws.Range("H2").Value = 0.106963
$ws.Range("I2").Value = 0.002412342638581826
$ws.Range("J2").Value = 0.002412342638581825
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.211751333333333
$ws.Range("N2").Value = 9.635254
$ws.Range("O2").Value = 0.9274105493513785
$ws.Range("P2").Value = 0.9274105493513782
$ws.Range("Q2").Value = 0.1145128526224445
$ws.Range("R2").Value = 1.030615673602
$ws.Range("S2").Value = 0.002237232011670925
$ws.Range("T2").Value = 0.002237232011670924

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Reln"
$ws.Range("C3").Value = "Lrp8"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.03565433333333334
$ws.Range("H3").Value = 0.106963
$ws.Range("I3").Value = 0.002412342638581826
$ws.Range("J3").Value = 0.002412342638581825
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.2513873333333334
$ws.Range("N3").Value = 0.754162
$ws.Range("O3").Value = 0.07258945064862164
$ws.Range("P3").Value = 0.07258945064862163
$ws.Range("Q3").Value = 0.008963047778444446
$ws.Range("R3").Value = 0.080667430006
$ws.Range("S3").Value = 0.0001751106269109012
$ws.Range("T3").Value = 0.0001751106269109011

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Reln"
$ws.Range("C4").Value = "Lrp8"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 8.066615333333333
$ws.Range("H4").Value = 24.199846
$ws.Range("I4").Value = 0.5457805068380079
$ws.Range("J4").Value = 0.5457805068380079
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.211751333333333
$ws.Range("N4").Value = 9.635254
$ws.Range("O4").Value = 0.9274105493513785
$ws.Range("P4").Value = 0.9274105493513782
$ws.Range("Q4").Value = 25.90796255232044
$ws.Range("R4").Value = 233.171662970884
$ws.Range("S4").Value = 0.5061625996719107
$ws.Range("T4").Value = 0.5061625996719106

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Reln"
$ws.Range("C5").Value = "Lrp8"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 8.066615333333333
$ws.Range("H5").Value = 24.199846
$ws.Range("I5").Value = 0.5457805068380079
$ws.Range("J5").Value = 0.5457805068380079
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.2513873333333334
$ws.Range("N5").Value = 0.754162
$ws.Range("O5").Value = 0.07258945064862164
$ws.Range("P5").Value = 0.07258945064862163
$ws.Range("Q5").Value = 2.027844917672445
$ws.Range("R5").Value = 18.250604259052
$ws.Range("S5").Value = 0.03961790716609728
$ws.Range("T5").Value = 0.03961790716609728

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Reln"
$ws.Range("C6").Value = "Lrp8"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 6.677692666666666
$ws.Range("H6").Value = 20.033078
$ws.Range("I6").Value = 0.4518071505234102
$ws.Range("J6").Value = 0.4518071505234102
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.211751333333333
$ws.Range("N6").Value = 9.635254
$ws.Range("O6").Value = 0.9274105493513785
$ws.Range("P6").Value = 0.9274105493513782
$ws.Range("Q6").Value = 21.44708832575689
$ws.Range("R6").Value = 193.023794931812
$ws.Range("S6").Value = 0.4190107176677968
$ws.Range("T6").Value = 0.4190107176677967

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Reln"
$ws.Range("C7").Value = "Lrp8"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 6.677692666666666
$ws.Range("H7").Value = 20.033078
$ws.Range("I7").Value = 0.4518071505234102
$ws.Range("J7").Value = 0.4518071505234102
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.2513873333333334
$ws.Range("N7").Value = 0.754162
$ws.Range("O7").Value = 0.07258945064862164
$ws.Range("P7").Value = 0.07258945064862163
$ws.Range("Q7").Value = 1.678687352292889
$ws.Range("R7").Value = 15.108186170636
$ws.Range("S7").Value = 0.03279643285561346
$ws.Range("T7").Value = 0.03279643285561345
